$wb = $excel.ActiveWorkbook

# Go to Sheet 2 ("Number two") and write the new cell value
$ws2 = $wb.Worksheets.Item("Number two")
$ws2.Range("B2").Value = "Sheet 2, cell 2 B"

# Select B3 on sheet 2 and activate the sheet, making it the active tab
$ws2.Activate()
$ws2.Range("B3").Select()
